$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Best Bound" / "Best Objective" header labels between C3 and D3
$ws.Range("C3").Value = "Best Bound"
$ws.Range("D3").Value = "Best Objective"

# Move the active selection to G4 (matches recorded selection in the saved file)
$ws.Range("G4").Select()
